# Generate Report for Handoff
#
# The "66c1a4d9-56e5-419e-9cb6-fe0ca0d55d46.md" file was handed off again for
# translation, so its "Latest Handoff Date"/"Latest Handoff Datetime" values
# (row 5 in every sheet) need to be refreshed to the new handoff timestamps.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest Handoff Date column (D) for row 5
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-32-19 08:32:12"

# "zh-cn" sheet: Latest Handoff Datetime column (E) for row 5
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-19 08:32:09"

# "de-de" sheet: Latest Handoff Datetime column (E) for row 5
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-19 08:32:12"
